$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Remove the J1 (C64-EXPANSION-PORT) row - Excel shifts all following rows up by one
$ws.Rows(2).Delete()

# Note the hardware work-around on the C1, C2 electrolytic capacitor row
# (this row is now row 3 after the delete above)
$ws.Range("E3").Value = "Used 22uF"

# Restore the view state: scrolled/zoomed in on the notes column
[void]$ws.Range("E4").Select()
$excel.ActiveWindow.Zoom = 160
